$wb = $excel.ActiveWorkbook

# Update the "Veda" sheet: fuel lookup table row 10 in column F changes
# from "wind" to "windon" (this intentionally breaks the VLOOKUP in D7,
# whose lookup key B7 is still "wind").
$wsVeda = $wb.Worksheets.Item("Veda")
$wsVeda.Range("F10").Value = "windon"

# Update the "historical_data_long" sheet: every "wind" entry in column A
# (the model_fuel column) is renamed to "windon".
$wsHist = $wb.Worksheets.Item("historical_data_long")
$lastRow = $wsHist.UsedRange.Rows.Count
for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $wsHist.Cells.Item($r, 1)
    if ($cell.Value2 -eq "wind") {
        $cell.Value = "windon"
    }
}
